$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BASE_FICHAS")
Write-Host $ws.Name
